$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.426.16"
$ws.Range("E2").Value = "  +3.27%  "

$ws.Range("D3").Value = "3.257.52"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'211.78"
$ws.Range("E5").Value = "  -3.18%  "

$ws.Range("D6").Value = "'625.11"
$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("E7").Value = "  +11.81%  "

$ws.Range("D8").Value = "'0.685"
$ws.Range("E8").Value = "  +15.59%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "3.255.33"
$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").Value = "'0.573"
$ws.Range("E11").Value = "  -5.12%  "

$ws.Range("E12").Value = "  +7.68%  "

$ws.Range("E13").Value = "  -9.09%  "

$ws.Range("D14").Value = "3.858.02"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").Value = "'33.61"
$ws.Range("E15").Value = "  -2.52%  "

$ws.Range("D16").Value = "'5.30"
$ws.Range("E16").Value = "  -2.43%  "

$ws.Range("D17").Value = "87.116.67"
$ws.Range("E17").Value = "  +3.24%  "

$ws.Range("D18").Value = "3.261.11"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").Value = "'3.12"
$ws.Range("E19").Value = "  -2.52%  "

$ws.Range("D20").Value = "'13.96"
$ws.Range("E20").Value = "  -4.05%  "

$ws.Range("D21").Value = "'432.05"
$ws.Range("E21").Value = "  -4.16%  "

$ws.Range("D22").Value = "'8.83"
$ws.Range("E22").Value = "  -4.27%  "

$ws.Range("D23").Value = "'5.29"
$ws.Range("E23").Value = "  +0.96%  "

$ws.Range("D24").Value = "'7.26"
$ws.Range("E24").Value = "  -2.82%  "

$ws.Range("D25").Value = "'12.42"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("D26").Value = "'5.08"
$ws.Range("E26").Value = "  -5.80%  "

$ws.Range("D27").Value = "3.412.34"
$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("D28").Value = "'75.88"
$ws.Range("E28").Value = "  -3.22%  "

$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  +10.41%  "

$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").Value = "'8.71"
$ws.Range("E33").Value = "  -5.94%  "

$ws.Range("D34").Value = "'542.69"
$ws.Range("E34").Value = "  -6.04%  "

$ws.Range("E35").Value = "  -9.24%  "

$ws.Range("E36").Value = "  -4.27%  "

$ws.Range("D37").Value = "'6.91"
$ws.Range("E37").Value = "  +9.24%  "

$ws.Range("E38").Value = "  -11.42%  "

$ws.Range("D39").Value = "'22.33"

$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").Value = "'21.68"

$ws.Range("D42").Value = "'0.391"
$ws.Range("E42").Value = "  -5.44%  "

$ws.Range("D43").Value = "'1.99"
$ws.Range("E43").Value = "  -3.57%  "

$ws.Range("E44").Value = "  -5.85%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "'154.84"
$ws.Range("E46").Value = "  -2.76%  "

$ws.Range("D47").Value = "'178.81"
$ws.Range("E47").Value = "  -6.69%  "

$ws.Range("D48").Value = "'44.87"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("E49").Value = "  -4.75%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.124"
$ws.Range("E50").Value = "  +12.99%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.19"
$ws.Range("E51").Value = "  -1.35%  "
